$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.107.77'
$ws.Range("E2").Value = '  +4.85%  '

$ws.Range("D3").Value = '2.610.94'
$ws.Range("E3").Value = '  +5.04%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.31'
$ws.Range("E5").Value = '  +3.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.12'
$ws.Range("E6").Value = '  +3.08%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  +1.93%  '

$ws.Range("D9").Value = '2.610.21'
$ws.Range("E9").Value = '  +4.96%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +16.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.164'
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  +4.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.05'
$ws.Range("E13").Value = '  +1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.85'
$ws.Range("E15").Value = '  +5.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000184'
$ws.Range("E16").Value = '  +7.98%  '

$ws.Range("D17").Value = '71.091.03'
$ws.Range("E17").Value = '  +4.92%  '

$ws.Range("D18").Value = '2.609.44'
$ws.Range("E18").Value = '  +6.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.95'
$ws.Range("E19").Value = '  +7.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '379.32'
$ws.Range("E20").Value = '  +8.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.45'
$ws.Range("E21").Value = '  +5.48%  '

$ws.Range("E22").Value = '  +1.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.10'
$ws.Range("E23").Value = '  +1.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.42'
$ws.Range("E24").Value = '  +5.03%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.88'
$ws.Range("E26").Value = '  +11.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("E27").Value = '  +9.63%  '

$ws.Range("D28").Value = '2.747.28'
$ws.Range("E28").Value = '  +6.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("D30").Value = '0.0₃0953'
$ws.Range("E30").Value = '  +6.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '531.66'
$ws.Range("E31").Value = '  +6.80%  '

$ws.Range("E32").Value = '  +4.82%  '

$ws.Range("E33").Value = '  +6.64%  '

$ws.Range("E34").Value = '  +4.51%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.94'
$ws.Range("E36").Value = '  +0.54%  '

$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.21'
$ws.Range("E38").Value = '  +4.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.95'
$ws.Range("E39").Value = '  +1.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.38'
$ws.Range("E40").Value = '  +6.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.85'
$ws.Range("E41").Value = '  +6.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.09'
$ws.Range("E42").Value = '  +5.97%  '

$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("E44").Value = '  +5.62%  '

$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.57'
$ws.Range("E46").Value = '  +4.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.37'
$ws.Range("E47").Value = '  +4.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.67'
$ws.Range("E48").Value = '  +4.13%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0269'
$ws.Range("E49").Value = '  +5.67%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.533'
$ws.Range("E50").Value = '  +4.38%  '

$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  +6.93%  '

Write-Host "Applied updates"